$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.049.87'
$ws.Range("E2").Value = '  +0.53%  '

$ws.Range("D3").Value = '2.756.84'
$ws.Range("E3").Value = '  +0.60%  '

$ws.Range("E4").Value = '  +0.13%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '576.62'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.41%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '159.06'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.09%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.27%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.603'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.60%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.110'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.64%  '

$ws.Range("B10").Value = 'TRON'
$ws.Range("C10").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.165'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.84%  '

$ws.Range("B11").Value = 'Toncoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.83'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -14.40%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.388'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.85%  '

$ws.Range("D13").Value = '3.248.17'
$ws.Range("E13").Value = '  +0.57%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.96'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.66%  '

$ws.Range("D15").Value = '63.682.81'
$ws.Range("E15").Value = '  +0.01%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000152'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.80%  '

$ws.Range("D17").Value = '2.763.07'
$ws.Range("E17").Value = '  +0.32%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.21'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.46%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.87'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.69%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '358.02'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.74%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.77'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.70%  '

$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.00'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.60%  '

$ws.Range("B23").Value = 'Polygon'
$ws.Range("C23").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.536'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.78%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.58'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.60%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.171'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.43%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.61'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.30%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.998'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.15%  '

$ws.Range("D28").Value = '0.0₃0912'
$ws.Range("E28").Value = '  -0.72%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.34'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.03%  '

$ws.Range("E30").Value = '  -3.52%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.24'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.06%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '168.83'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.78%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.94'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.39%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '20.21'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.20%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.49'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.58%  '

$ws.Range("E36").Value = '  +0.09%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.81'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.82%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.24%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '350.53'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.60%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.34'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.78%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.20'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.38%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '39.19'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.12%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '21.57'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.38%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '21.91'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.64%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0591'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.45%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '137.27'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.13%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0256'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.53%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.633'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.97%  '

$ws.Range("E49").Value = '  -0.52%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.998'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.16%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '11.03'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.14%  '
